$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all existing data rows
# (rows 2-516) from 2023-09-13 (45182) to 2023-09-15 (45184).
$ws.Range("C2:C516").Value = 45184

# Row 516 picks up an explicit row height in the target file.
$ws.Rows.Item(516).RowHeight = 15

# Append the new record as row 517.
$ws.Range("A517").Value = "A 43230-2023"

$ws.Range("B517").Value = 45183
$ws.Range("B517").NumberFormat = "YYYY-MM-DD"

$ws.Range("C517").Value = 45184
$ws.Range("C517").NumberFormat = "YYYY-MM-DD"

$ws.Range("D517").Value = "UPPSALA LÄN"
$ws.Range("E517").Value = "ENKÖPING"

$ws.Range("G517").Value = 2.1
$ws.Range("H517").Value = 0
$ws.Range("I517").Value = 0
$ws.Range("J517").Value = 0
$ws.Range("K517").Value = 0
$ws.Range("L517").Value = 0
$ws.Range("M517").Value = 0
$ws.Range("N517").Value = 0
$ws.Range("O517").Value = 0
$ws.Range("P517").Value = 0
$ws.Range("Q517").Value = 0

$ws.Range("R517").WrapText = $true
